$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.583.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.57%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.913.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.72%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.74%  "

# Row 6
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.700"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.00%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "45.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.88%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.371"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.55%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.26%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0762"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.84%  "

# Row 12
$ws.Range("E12").Value = "  +2.51%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.38%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.811"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.57%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.191.63"
$ws.Range("D15").Style = "Normal"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.75%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.918.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.98%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.606.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.72%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.29%  "

# Row 20
$ws.Range("E20").Value = "  +4.69%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "250.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.22%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.07%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.22%  "

# Row 24
$ws.Range("E24").Value = "  -2.99%  "

# Row 25
$ws.Range("E25").Value = "  -0.08%  "

# Row 26
$ws.Range("E26").Value = "  +2.06%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "168.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.25%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.17%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.32%  "

# Row 30
$ws.Range("E30").Value = "  +1.64%  "

# Row 31
$ws.Range("E31").Value = "  +6.28%  "

# Row 32
$ws.Range("E32").Value = "  +4.60%  "

# Row 33
$ws.Range("E33").Value = "  +3.97%  "

# Row 34
$ws.Range("E34").Value = "  +23.69%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.91"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.62%  "

# Row 36
$ws.Range("E36").Value = "  +0.09%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.52"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.38%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.881"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.44%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +50.88%  "

# Row 40
$ws.Range("E40").Value = "  +4.13%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.77%  "

# Row 42
$ws.Range("E42").Value = "  +3.91%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.53"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.00%  "

# Row 44
$ws.Range("E44").Value = "  +23.23%  "

# Row 45
$ws.Range("E45").Value = "  +3.16%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.343.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.70%  "

# Row 47
$ws.Range("E47").Value = "  -0.72%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0812"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.66%  "

# Row 49
$ws.Range("E49").Value = "  +2.26%  "

# Row 50
$ws.Range("E50").Value = "  +1.52%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.52%  "
